# Day-3 topic reorder: "GitOps" and "Docker" swap places in the title,
# and the "GitOps Grundlagen" sub-bullet moves from the first to the
# last position of the Day-3 bullet block. Applied identically to the
# two (near-duplicate) Agenda slides. Also refreshes the day-off-by-one
# date stamp baked into two of the slide masters.

$p = $ppt.ActivePresentation

function Set-Day3Agenda($slideIndex) {
    $s = $p.Slides.Item($slideIndex)
    $shp = $s.Shapes.Item(1)
    $tr = $shp.TextFrame.TextRange

    # --- Title paragraph (#13): "Tag 3 - GitOps, Docker, Deployment-Strategien"
    #     becomes "Tag 3 - Docker, GitOps, Deployment-Strategien".
    #     Runs: [Tag 3 - ][GitOps][, Docker, ][Deployment][-Strategien]
    #     Only runs 1 and 3 change text; runs 2/4/5 stay untouched so
    #     their run-level formatting (err="1" spell markers) survives.
    $dash = [char]0x2013
    $titlePara = $tr.Paragraphs(13, 1)
    $base = $titlePara.Start
    $off = 0

    $run1Old = 8   # "Tag 3 - " (en dash)
    $run1New = "Tag 3 $dash Docker, "
    $c = $tr.Characters($base + $off, $run1Old)
    $c.Text = $run1New
    $off += $run1New.Length

    $off += 6      # run 2 "GitOps" - unchanged, just skip over it

    $run3Old = 10  # ", Docker, "
    $run3New = ", "
    $c = $tr.Characters($base + $off, $run3Old)
    $c.Text = $run3New
    $off += $run3New.Length

    # remaining runs (Deployment / -Strategien) untouched

    # --- Bullet block (#14-#17): rotate so that paragraph #14's content
    #     ("GitOps Grundlagen") ends up in slot #17, with #15/#16/#17
    #     shifting up to #14/#15/#16. Run layouts line up 1:1 between
    #     source and destination slot, so we can rewrite run-by-run and
    #     keep each run's own formatting (e.g. err="1" on "GitOps" /
    #     "Tagged") intact.

    # Snapshot of the paragraphs' run texts *before* any edits, in the
    # original (pre-change) order.
    $p14Runs = @("GitOps", " Grundlagen")
    $p15Runs = @("Entwicklung mit Docker")
    $p16Runs = @("Container/Docker-Registry")
    $p17Runs = @("Erstellen von Release- und ", "Tagged", "-Images")

    # New content for each paragraph slot, taken from the rotation.
    $newForP14 = $p15Runs
    $newForP15 = $p16Runs
    $newForP16 = $p17Runs
    $newForP17 = $p14Runs

    $oldLensForP14 = @($p14Runs[0].Length, $p14Runs[1].Length)
    $oldLensForP15 = @($p15Runs[0].Length)
    $oldLensForP16 = @($p16Runs[0].Length)
    $oldLensForP17 = @($p17Runs[0].Length, $p17Runs[1].Length, $p17Runs[2].Length)

    Set-ParagraphRuns $tr 14 $oldLensForP14 $newForP14
    Set-ParagraphRuns $tr 15 $oldLensForP15 $newForP15
    Set-ParagraphRuns $tr 16 $oldLensForP16 $newForP16
    Set-ParagraphRuns $tr 17 $oldLensForP17 $newForP17
}

function Set-ParagraphRuns($tr, $paraIndex, $oldRunLens, $newRunTexts) {
    $para = $tr.Paragraphs($paraIndex, 1)
    $base = $para.Start
    $off = 0
    for ($i = 0; $i -lt $oldRunLens.Count; $i++) {
        $oldLen = $oldRunLens[$i]
        $newText = $newRunTexts[$i]
        $c = $tr.Characters($base + $off, $oldLen)
        $c.Text = $newText
        $off += $newText.Length
    }
}

Set-Day3Agenda 2
Set-Day3Agenda 3

# --- Slide-master date stamps: "18.06.2024" -> "19.06.2024".
# (These live in a plain rectangle shape containing an auto-date field;
# rewriting its text is the only way to reach it through the object
# model, which bakes the field down to literal text - matching what
# PowerPoint itself does when the cached field text is refreshed.)
function Set-MasterDate($designIndex) {
    $design = $p.Designs.Item($designIndex)
    $master = $design.SlideMaster
    for ($i = 1; $i -le $master.Shapes.Count; $i++) {
        $shp = $master.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.HasText) {
                $t = $shp.TextFrame.TextRange.Text
                if ($t -eq "18.06.2024") {
                    $shp.TextFrame.TextRange.Text = "19.06.2024"
                }
            }
        }
    }
}

for ($di = 1; $di -le $p.Designs.Count; $di++) {
    Set-MasterDate $di
}
